$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student Scores")

# Update the "Written Assignment (15)" column (D) values for a few students
$ws.Range("D4").Value = 8
$ws.Range("D11").Value = 8
$ws.Range("D12").Value = 8

# Update the active selection on the sheet from D13 to I6
$ws.Range("I6").Select()
